$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "General" cells that would otherwise be mis-parsed as numbers to stay text
# (values such as "1.009" or "11.10" must remain literal text, matching the source data)

$ws.Range("D2").Value = '28.507.35'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '1.876.06'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("E4").Value = '  -0.87%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.93'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5082'
$ws.Range("E7").Value = '  -1.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3909'
$ws.Range("E8").Value = '  -0.72%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08378'
$ws.Range("E9").Value = '  +0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.228'
$ws.Range("D12").Value = '1.872.58'
$ws.Range("E12").Value = '  +0.39%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.41'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.257'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.011'
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001105'
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.34'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06733'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.928'
$ws.Range("E21").Value = '  -1.17%  '
$ws.Range("D22").Value = '28.515.08'
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.11'
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.234'
$ws.Range("E24").Value = '  -1.06%  '
$ws.Range("D25").Value = '2.090.40'
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.71'
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.62'
$ws.Range("E27").Value = '  -0.77%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.392'
$ws.Range("E28").Value = '  -0.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '125.93'
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.041'
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.773'
$ws.Range("E32").Value = '  -1.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.614'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.02454'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06561'
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.825'
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.065'
$ws.Range("E38").Value = '  +1.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.256'
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.193'
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6401'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.10'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6026'
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.08'
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.688'
$ws.Range("E46").Value = '  -0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.007'
$ws.Range("E47").Value = '  -0.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.217'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '121.96'
$ws.Range("E49").Value = '  +0.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.130'
$ws.Range("E50").Value = '  -12.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06849'
$ws.Range("E51").Value = '  -0.76%  '
